# Update sensor distance columns E (x_m) and F (y_m) with recalculated values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(2, 21.00858044937448, 5.328837811611383),
    @(3, 21.00858044937448, 5.328837811611383),
    @(4, 20.27513380211692, 1.670248450302496),
    @(5, 20.27513380211692, 1.670248450302496),
    @(6, 26.6628488107649, 1.5946562518322),
    @(7, 26.6628488107649, 1.5946562518322),
    @(8, 26.6628488107649, 1.5946562518322),
    @(9, 26.6628488107649, 1.5946562518322),
    @(10, 0.3685083442959955, 1.635798685098902),
    @(11, 0.3685083442959955, 1.635798685098902),
    @(12, 0.3685083442959955, 1.635798685098902),
    @(13, 0.3685083442959955, 1.635798685098902),
    @(14, 10.10647680343427, 6.057379903310235),
    @(15, 10.10647680343427, 6.057379903310235),
    @(16, 10.10647680343427, 6.057379903310235),
    @(17, 10.10647680343427, 6.057379903310235),
    @(18, 12.01411921803544, 4.705269775963187),
    @(19, 12.01411921803544, 4.705269775963187),
    @(20, 12.01411921803544, 4.705269775963187),
    @(21, 12.01411921803544, 4.705269775963187),
    @(22, 11.06029801073486, 5.381324839636711),
    @(23, 11.06029801073486, 5.381324839636711),
    @(24, 11.06029801073486, 5.381324839636711),
    @(25, 11.06029801073486, 5.381324839636711),
    @(26, 11.06029801073486, 5.381324839636711),
    @(27, 11.06029801073486, 5.381324839636711),
    @(28, 13.92176163263661, 6.057379903310235),
    @(29, 13.92176163263661, 6.057379903310235),
    @(30, 13.92176163263661, 6.057379903310235),
    @(31, 13.92176163263661, 6.057379903310235),
    @(32, 15.82940404723777, 4.705269775963187),
    @(33, 15.82940404723777, 4.705269775963187),
    @(34, 15.82940404723777, 4.705269775963187),
    @(35, 15.82940404723777, 4.705269775963187),
    @(36, 14.87558283993719, 5.381324839636711),
    @(37, 14.87558283993719, 5.381324839636711),
    @(38, 14.87558283993719, 5.381324839636711),
    @(39, 14.87558283993719, 5.381324839636711),
    @(40, 14.87558283993719, 5.381324839636711),
    @(41, 14.87558283993719, 5.381324839636711),
    @(42, 17.73704646183894, 6.057379903310235),
    @(43, 17.73704646183894, 6.057379903310235),
    @(44, 17.73704646183894, 6.057379903310235),
    @(45, 17.73704646183894, 6.057379903310235),
    @(46, 19.64468887644011, 4.705269775963187),
    @(47, 19.64468887644011, 4.705269775963187),
    @(48, 19.64468887644011, 4.705269775963187),
    @(49, 19.64468887644011, 4.705269775963187),
    @(50, 18.69086766913952, 5.381324839636711),
    @(51, 18.69086766913952, 5.381324839636711),
    @(52, 18.69086766913952, 5.381324839636711),
    @(53, 18.69086766913952, 5.381324839636711),
    @(54, 18.69086766913952, 5.381324839636711),
    @(55, 18.69086766913952, 5.381324839636711)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
}
